# Update "想去人数" (F column) counts for a few events on both the
# "展览" sheet and the "全部类型" sheet, reflecting freshly re-generated
# data output (gh-pages build at 456a3b4).

$wb = $excel.ActiveWorkbook

# Sheet "展览" (index 1)
$wsExhibit = $wb.Worksheets.Item(1)
$wsExhibit.Range("F3").Value = 575
$wsExhibit.Range("F7").Value = 52
$wsExhibit.Range("F11").Value = 4684
$wsExhibit.Range("F12").Value = 4477

# Sheet "全部类型" (index 4)
$wsAll = $wb.Worksheets.Item(4)
$wsAll.Range("F3").Value = 576
$wsAll.Range("F7").Value = 52
$wsAll.Range("F11").Value = 4684
$wsAll.Range("F12").Value = 4477
